# Word COM-interop script implementing the commit "Added new API examples".
#
# Semantics: in the "GroupBy(Selector)" template example paragraph, the
# outer "<<foreach [m in ...]>> ... <</foreach>>" block is tagged "#header1"
# and the inner "<<foreach [c in m]>> ... <</foreach>>" block is tagged
# "#header2" -- i.e. the opening/closing tag markers gain " #header1" /
# " #header2" suffixes just before their closing ">>".

$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive, unique leading text.
$full = $d.Content.Text
$pstart = $full.IndexOf("Managers.GroupBy(p")
if ($pstart -lt 0) { throw "anchor paragraph not found" }

# --- Step 1: rewrite the three tag-delimiter runs (text-content only) ---
# These three Range.Text assignments turn:
#   ...})]>><<foreach [c in m]>>Age = ...; <</foreach>><</foreach>>
# into:
#   ...})] #header1>><<foreach [c in m] #header2>>Age = ...; <</foreach #header2>><</foreach #header1>>
# Applied right-to-left so earlier offsets stay valid while later ones shift.

# "<</foreach>><</foreach>>" (rel 107..131) -> "<</foreach #header2>><</foreach #header1>>"
$r3 = $d.Range($pstart + 107, $pstart + 131)
if ($r3.Text -ne "<</foreach>><</foreach>>") { throw "hunk3 anchor mismatch: $($r3.Text)" }
$r3.Text = "<</foreach #header2>><</foreach #header1>>"

# "]>>" (rel 64..67, the inner foreach closer) -> "] #header2>>"
$r2 = $d.Range($pstart + 64, $pstart + 67)
if ($r2.Text -ne "]>>") { throw "hunk2 anchor mismatch: $($r2.Text)" }
$r2.Text = "] #header2>>"

# "]>><<" (rel 44..49, the outer foreach closer / inner foreach opener) -> "] #header1>><<"
$r1 = $d.Range($pstart + 44, $pstart + 49)
if ($r1.Text -ne "]>><<") { throw "hunk1 anchor mismatch: $($r1.Text)" }
$r1.Text = "] #header1>><<"

# --- Step 2: re-split the resulting run back into one <w:r> per literal/tag ---
# The three Range.Text writes above coalesce every same-formatted run from
# " })" to the end of the paragraph into a single run. Re-impose the exact
# run boundaries (original ones plus the new ones introduced by the commit)
# by toggling a character-level Font property on/off over each sub-range:
# Word always materialises an explicit run split at a Range used for a
# direct-formatting write, and restoring the same value leaves the visible
# formatting untouched.
$full2 = $d.Content.Text
$pstart2 = $full2.IndexOf("Managers.GroupBy(p")
$tailStart = $pstart2 + 41

$b = $d.Range($tailStart + 0, $tailStart + 3)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 3, $tailStart + 4)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 4, $tailStart + 13)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 13, $tailStart + 17)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 17, $tailStart + 25)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 25, $tailStart + 26)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 26, $tailStart + 32)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 32, $tailStart + 33)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 33, $tailStart + 34)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 34, $tailStart + 41)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 41, $tailStart + 42)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 42, $tailStart + 44)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 44, $tailStart + 50)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 50, $tailStart + 52)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 52, $tailStart + 53)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 53, $tailStart + 58)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 58, $tailStart + 59)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 59, $tailStart + 61)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 61, $tailStart + 62)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 62, $tailStart + 70)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 70, $tailStart + 72)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 72, $tailStart + 73)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 73, $tailStart + 79)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 79, $tailStart + 80)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 80, $tailStart + 82)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 82, $tailStart + 83)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 83, $tailStart + 84)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 84, $tailStart + 94)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 94, $tailStart + 95)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 95, $tailStart + 102)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 102, $tailStart + 103)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 103, $tailStart + 105)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 105, $tailStart + 115)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 115, $tailStart + 116)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 116, $tailStart + 124)
$b.Font.Bold = 1
$b.Font.Bold = 0
$b = $d.Range($tailStart + 124, $tailStart + 126)
$b.Font.Bold = 1
$b.Font.Bold = 0

# --- Step 3: sanity-check the resulting text ---
$full3 = $d.Content.Text
$pstart3 = $full3.IndexOf("Managers.GroupBy(p")
$pend3 = $pstart3 + 41 + 126
$result = $full3.Substring($pstart3, $pend3 - $pstart3)
$expected = "Managers.GroupBy(p => new { p.Age, p.Name })] #header1>><<foreach [c in m] #header2>>Age = <<[c.Age]>>, Name = <<[c.Name]>>; <</foreach #header2>><</foreach #header1>>"
Write-Output $result
if ($result -ne $expected) { throw "post-condition mismatch" }
